# Add a new column M (year 2021) to the right of the existing L column (year 2020),
# then switch the numeric data range from the custom "0.0" number format to
# the default "General" format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column M: copy the formatting from column L (the previous last
#     data column) onto M for every row that currently carries formatting,
#     then fill in the 2021 figures. ---

# Row 2 (thin divider row under the header)
$ws.Cells.Item(2, 12).Copy()
$ws.Cells.Item(2, 13).PasteSpecial(-4122)   # xlPasteFormats

# Row 3 (year headers) - M3 should look like F3 (no-fill variant), matching
# the alternating fill pattern already used across the year columns.
$ws.Cells.Item(3, 6).Copy()
$ws.Cells.Item(3, 13).PasteSpecial(-4122)
$ws.Cells.Item(3, 13).Value = 2021

# Row 4 (256 kbit/s - 2Mbit/s)
$ws.Cells.Item(4, 12).Copy()
$ws.Cells.Item(4, 13).PasteSpecial(-4122)
$ws.Cells.Item(4, 13).Value = 7105

# Row 5 (2 Mbit/s - 10 Mbit/s)
$ws.Cells.Item(5, 12).Copy()
$ws.Cells.Item(5, 13).PasteSpecial(-4122)
$ws.Cells.Item(5, 13).Value = 81079

# Row 6 (>= 10 Mbit/s)
$ws.Cells.Item(6, 12).Copy()
$ws.Cells.Item(6, 13).PasteSpecial(-4122)
$ws.Cells.Item(6, 13).Value = 214139

# --- Switch the data block (now D4:M6) from the custom "0.0" number
#     format to the plain "General" format. ---
$ws.Range("D4:M6").NumberFormat = "general"

# Clear the lingering selection/active-cell marker left on C2 so the sheet
# reopens with the default A1 selection.
[void]$ws.Range("A1").Select()
